# Weekly "cryptos" data refresh (GitHub Actions scheduled update).
# Refreshes the Price (D) and Volume(1h) (E) columns for every listed coin,
# and reflects the ranking swap where NEARProtocol (rows 49) now edges out
# Cronos (row 50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prices are stored as plain text (e.g. "30.647.78" uses "." as a thousands
# separator, and small caps like "0.9991" would otherwise be auto-parsed into
# a number by Excel). Mirror what a user gets by typing a leading apostrophe
# ('234.43) whenever the new price text looks numeric, so the cell keeps its
# original text type; leave alphanumeric-formatted prices (with multiple dots)
# untouched since Excel already treats those as text.
function Set-PriceText($cellRef, $text) {
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $ws.Range($cellRef).Value = "'" + $text
    } else {
        $ws.Range($cellRef).Value = $text
    }
}

Set-PriceText "D2" '30.698.30'
$ws.Range("E2").Value = '  +1.24%  '

Set-PriceText "D3" '1.865.86'
$ws.Range("E3").Value = '  +0.22%  '

Set-PriceText "D4" '0.9982'
$ws.Range("E4").Value = '  -0.18%  '

Set-PriceText "D5" '234.43'
$ws.Range("E5").Value = '  +0.25%  '

Set-PriceText "D6" '0.9987'
$ws.Range("E6").Value = '  -0.09%  '

Set-PriceText "D7" '0.4714'
$ws.Range("E7").Value = '  -0.93%  '

Set-PriceText "D8" '0.2759'
$ws.Range("E8").Value = '  +0.22%  '

Set-PriceText "D9" '0.06368'
$ws.Range("E9").Value = '  -1.20%  '

Set-PriceText "D10" '17.46'
$ws.Range("E10").Value = '  +8.23%  '

Set-PriceText "D11" '1.861.40'
$ws.Range("E11").Value = '  -0.10%  '

Set-PriceText "D12" '0.07475'
$ws.Range("E12").Value = '  +0.62%  '

Set-PriceText "D13" '4.974'
$ws.Range("E13").Value = '  -0.41%  '

Set-PriceText "D14" '85.17'
$ws.Range("E14").Value = '  -1.02%  '

Set-PriceText "D15" '0.6312'
$ws.Range("E15").Value = '  -0.38%  '

Set-PriceText "D16" '30.624.15'
$ws.Range("E16").Value = '  +1.04%  '

Set-PriceText "D17" '241.58'
$ws.Range("E17").Value = '  +3.86%  '

Set-PriceText "D18" '0.9991'
$ws.Range("E18").Value = '  -0.06%  '

Set-PriceText "D19" '12.73'
$ws.Range("E19").Value = '  -0.73%  '

Set-PriceText "D20" '0.000007382'
$ws.Range("E20").Value = '  -0.04%  '

Set-PriceText "D21" '0.9981'
$ws.Range("E21").Value = '  -0.18%  '

Set-PriceText "D22" '4.988'
$ws.Range("E22").Value = '  -2.29%  '

Set-PriceText "D23" '5.963'
$ws.Range("E23").Value = '  -0.99%  '

Set-PriceText "D24" '9.296'
$ws.Range("E24").Value = '  -0.05%  '

Set-PriceText "D25" '166.55'
$ws.Range("E25").Value = '  -0.45%  '

Set-PriceText "D26" '18.18'
$ws.Range("E26").Value = '  +1.36%  '

Set-PriceText "D27" '1.886'
$ws.Range("E27").Value = '  +1.32%  '

Set-PriceText "D28" '0.1029'
$ws.Range("E28").Value = '  +2.38%  '

Set-PriceText "D29" '1.377'
$ws.Range("E29").Value = '  -0.50%  '

Set-PriceText "D30" '4.113'
$ws.Range("E30").Value = '  -2.84%  '

Set-PriceText "D31" '3.862'
$ws.Range("E31").Value = '  -1.43%  '

Set-PriceText "D32" '0.04941'
$ws.Range("E32").Value = '  +0.80%  '

Set-PriceText "D33" '1.153'
$ws.Range("E33").Value = '  +0.28%  '

Set-PriceText "D34" '0.7093'
$ws.Range("E34").Value = '  -2.26%  '

Set-PriceText "D35" '2.698'
$ws.Range("E35").Value = '  +0.22%  '

Set-PriceText "D36" '0.01916'
$ws.Range("E36").Value = '  -1.25%  '

Set-PriceText "D37" '2.688'
$ws.Range("E37").Value = '  +2.08%  '

Set-PriceText "D38" '0.8845'
$ws.Range("E38").Value = '  -2.71%  '

Set-PriceText "D39" '1.999'
$ws.Range("E39").Value = '  +0.40%  '

Set-PriceText "D40" '106.03'
$ws.Range("E40").Value = '  +0.41%  '

Set-PriceText "D41" '0.9989'
$ws.Range("E41").Value = '  -0.09%  '

Set-PriceText "D42" '0.4104'
$ws.Range("E42").Value = '  -0.29%  '

Set-PriceText "D43" '5.556'
$ws.Range("E43").Value = '  +0.12%  '

Set-PriceText "D44" '7.208'
$ws.Range("E44").Value = '  +1.92%  '

Set-PriceText "D45" '0.1234'
$ws.Range("E45").Value = '  +2.13%  '

Set-PriceText "D46" '61.95'
$ws.Range("E46").Value = '  +0.98%  '

Set-PriceText "D47" '8.647'
$ws.Range("E47").Value = '  -1.39%  '

$ws.Range("E48").Value = '  +1.79%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-PriceText "D49" '1.383'
$ws.Range("E49").Value = '  -1.34%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-PriceText "D50" '0.05565'
$ws.Range("E50").Value = '  -0.78%  '

Set-PriceText "D51" '0.3713'
$ws.Range("E51").Value = '  -0.08%  '
